$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 123
$ws1.Range("F3").Value = 14

# Sheet "全部类型" (All types) mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 123
$ws4.Range("F3").Value = 14
